$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (B and C) ---
$ws.Columns.Item(2).ColumnWidth = 46.833333333333336
$ws.Columns.Item(3).ColumnWidth = 87.16666666666667

# --- Wrap text styling ---
$ws.Range("B1:C1").WrapText = $true
$ws.Range("B2:C71").WrapText = $true

# --- Row heights (per-row, matches target ht values) ---
$ws.Rows.Item(1).RowHeight = 16
$ws.Rows.Item(2).RowHeight = 304
$ws.Rows.Item(3).RowHeight = 208
$ws.Rows.Item(4).RowHeight = 256
$ws.Rows.Item(5).RowHeight = 208
$ws.Rows.Item(6).RowHeight = 320
$ws.Rows.Item(7).RowHeight = 350
$ws.Rows.Item(8).RowHeight = 365
$ws.Rows.Item(9).RowHeight = 272
$ws.Rows.Item(10).RowHeight = 224
$ws.Rows.Item(11).RowHeight = 335
$ws.Rows.Item(12).RowHeight = 409.6
$ws.Rows.Item(13).RowHeight = 380
$ws.Rows.Item(14).RowHeight = 272
$ws.Rows.Item(15).RowHeight = 112
$ws.Rows.Item(16).RowHeight = 240
$ws.Rows.Item(17).RowHeight = 208
$ws.Rows.Item(18).RowHeight = 272
$ws.Rows.Item(19).RowHeight = 272
$ws.Rows.Item(20).RowHeight = 335
$ws.Rows.Item(21).RowHeight = 380
$ws.Rows.Item(22).RowHeight = 350
$ws.Rows.Item(23).RowHeight = 272
$ws.Rows.Item(24).RowHeight = 208
$ws.Rows.Item(25).RowHeight = 272
$ws.Rows.Item(26).RowHeight = 409.6
$ws.Rows.Item(27).RowHeight = 365
$ws.Rows.Item(28).RowHeight = 272
$ws.Rows.Item(29).RowHeight = 32
$ws.Rows.Item(30).RowHeight = 224
$ws.Rows.Item(31).RowHeight = 208
$ws.Rows.Item(32).RowHeight = 224
$ws.Rows.Item(33).RowHeight = 208
$ws.Rows.Item(34).RowHeight = 320
$ws.Rows.Item(35).RowHeight = 350
$ws.Rows.Item(36).RowHeight = 335
$ws.Rows.Item(37).RowHeight = 272
$ws.Rows.Item(38).RowHeight = 208
$ws.Rows.Item(39).RowHeight = 272
$ws.Rows.Item(40).RowHeight = 380
$ws.Rows.Item(41).RowHeight = 365
$ws.Rows.Item(42).RowHeight = 272
$ws.Rows.Item(43).RowHeight = 32
$ws.Rows.Item(44).RowHeight = 208
$ws.Rows.Item(45).RowHeight = 192
$ws.Rows.Item(46).RowHeight = 224
$ws.Rows.Item(47).RowHeight = 256
$ws.Rows.Item(48).RowHeight = 288
$ws.Rows.Item(49).RowHeight = 350
$ws.Rows.Item(50).RowHeight = 320
$ws.Rows.Item(51).RowHeight = 272
$ws.Rows.Item(52).RowHeight = 240
$ws.Rows.Item(53).RowHeight = 272
$ws.Rows.Item(54).RowHeight = 409.6
$ws.Rows.Item(55).RowHeight = 350
$ws.Rows.Item(56).RowHeight = 272
$ws.Rows.Item(57).RowHeight = 128
$ws.Rows.Item(58).RowHeight = 240
$ws.Rows.Item(59).RowHeight = 192
$ws.Rows.Item(60).RowHeight = 256
$ws.Rows.Item(61).RowHeight = 272
$ws.Rows.Item(62).RowHeight = 335
$ws.Rows.Item(63).RowHeight = 365
$ws.Rows.Item(64).RowHeight = 365
$ws.Rows.Item(65).RowHeight = 272
$ws.Rows.Item(66).RowHeight = 288
$ws.Rows.Item(67).RowHeight = 272
$ws.Rows.Item(68).RowHeight = 409.6
$ws.Rows.Item(69).RowHeight = 409.6
$ws.Rows.Item(70).RowHeight = 272
$ws.Rows.Item(71).RowHeight = 32

# --- Column D relevance values ---
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(21, 4).Value = 1
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(23, 4).Value = 1
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(34, 4).Value = 1
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(38, 4).Value = 1
$ws.Cells.Item(39, 4).Value = 1
$ws.Cells.Item(40, 4).Value = 1
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(42, 4).Value = 1
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(45, 4).Value = 1
$ws.Cells.Item(46, 4).Value = 1
$ws.Cells.Item(47, 4).Value = 1
$ws.Cells.Item(48, 4).Value = 1
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(50, 4).Value = 1
$ws.Cells.Item(51, 4).Value = 1
$ws.Cells.Item(52, 4).Value = 1
$ws.Cells.Item(53, 4).Value = 1
$ws.Cells.Item(54, 4).Value = 1
$ws.Cells.Item(55, 4).Value = 1
$ws.Cells.Item(56, 4).Value = 1
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(59, 4).Value = 1
$ws.Cells.Item(60, 4).Value = 1
$ws.Cells.Item(61, 4).Value = 1
$ws.Cells.Item(62, 4).Value = 1
$ws.Cells.Item(63, 4).Value = 1
$ws.Cells.Item(64, 4).Value = 1
$ws.Cells.Item(65, 4).Value = 1
$ws.Cells.Item(66, 4).Value = 1
$ws.Cells.Item(67, 4).Value = 1
$ws.Cells.Item(68, 4).Value = 1
$ws.Cells.Item(69, 4).Value = 1
$ws.Cells.Item(70, 4).Value = 1
$ws.Cells.Item(71, 4).Value = 0

# --- View state: zoom + selection ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("D71").Select()

